$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, pushing existing rows 3-15 down to 4-16.
$ws.Rows.Item(3).Insert()

# Copy the date-cell style from row 2's date column into the newly inserted row's date cell
$ws.Cells.Item(2, 4).Copy()
$ws.Cells.Item(3, 4).PasteSpecial(-4122)  # xlPasteFormats

# Fill in the values for the new weekly entry (row 3)
$ws.Cells.Item(3, 1).Value = 10
$ws.Cells.Item(3, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(3, 3).Value = "La Araucanía"
$ws.Cells.Item(3, 4).Value = 44819
$ws.Cells.Item(3, 5).Value = 9
$ws.Cells.Item(3, 6).Value = 100112036
$ws.Cells.Item(3, 7).Value = "Caigua"
$ws.Cells.Item(3, 8).Value = "Sin especificar"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 100
$ws.Cells.Item(3, 11).Value = 20000
$ws.Cells.Item(3, 12).Value = 20000
$ws.Cells.Item(3, 13).Value = 20000
$ws.Cells.Item(3, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(3, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(3, 16).Value = 1333
$ws.Cells.Item(3, 17).Value = 15
$ws.Cells.Item(3, 18).Value = "Hortaliza"
